$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format before writing numeric-looking price strings,
# then restore Normal style so no stray cell-style ("s=") is left behind —
# matches the source diff, which only changes cell text, never styles.
$dCells = $ws.Range("D2:D51")
$dCells.NumberFormat = "@"

$ws.Range("D2").Value = "54.471.75"
$ws.Range("E2").Value = "  -5.61%  "
$ws.Range("D3").Value = "2.876.45"
$ws.Range("E3").Value = "  -9.05%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "468.70"
$ws.Range("E5").Value = "  -11.53%  "
$ws.Range("D6").Value = "125.85"
$ws.Range("E6").Value = "  -6.49%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "2.876.31"
$ws.Range("E8").Value = "  -9.04%  "
$ws.Range("D9").Value = "0.405"
$ws.Range("E9").Value = "  -10.33%  "
$ws.Range("E10").Value = "  -8.05%  "
$ws.Range("D11").Value = "0.0970"
$ws.Range("E11").Value = "  -13.27%  "
$ws.Range("E12").Value = "  -16.19%  "
$ws.Range("E13").Value = "  -4.04%  "
$ws.Range("D14").Value = "3.366.95"
$ws.Range("E14").Value = "  -9.11%  "
$ws.Range("D15").Value = "23.30"
$ws.Range("E15").Value = "  -10.08%  "
$ws.Range("D16").Value = "54.434.08"
$ws.Range("E16").Value = "  -5.59%  "
$ws.Range("D17").Value = "2.876.61"
$ws.Range("E17").Value = "  -8.99%  "
$ws.Range("E18").Value = "  -13.39%  "
$ws.Range("D19").Value = "5.34"
$ws.Range("E19").Value = "  -8.76%  "
$ws.Range("D20").Value = "11.42"
$ws.Range("E20").Value = "  -13.30%  "
$ws.Range("D21").Value = "7.05"
$ws.Range("E21").Value = "  -12.91%  "
$ws.Range("D22").Value = "300.22"
$ws.Range("E22").Value = "  -14.01%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  -10.00%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "0.0₃0797"
$ws.Range("E29").Value = "  -17.74%  "
$ws.Range("D30").Value = "6.13"
$ws.Range("E30").Value = "  -11.91%  "
$ws.Range("D31").Value = "6.19"
$ws.Range("E31").Value = "  -11.46%  "
$ws.Range("E32").Value = "  -9.89%  "
$ws.Range("D33").Value = "18.82"
$ws.Range("E33").Value = "  -13.15%  "
$ws.Range("E34").Value = "  -15.45%  "
$ws.Range("D35").Value = "141.57"
$ws.Range("E35").Value = "  -10.98%  "
$ws.Range("E36").Value = "  -16.36%  "
$ws.Range("D37").Value = "5.40"
$ws.Range("E37").Value = "  -13.88%  "
$ws.Range("D38").Value = "1.20"
$ws.Range("E38").Value = "  -15.08%  "
$ws.Range("D39").Value = "22.66"
$ws.Range("E39").Value = "  -14.21%  "
$ws.Range("D40").Value = "0.0617"
$ws.Range("E40").Value = "  -12.10%  "
$ws.Range("D41").Value = "2.902.74"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "35.16"
$ws.Range("E43").Value = "  -12.92%  "
$ws.Range("E46").Value = "  -14.15%  "
$ws.Range("E47").Value = "  -11.08%  "
$ws.Range("D48").Value = "2.028.66"
$ws.Range("E48").Value = "  -10.87%  "
$ws.Range("E49").Value = "  -14.65%  "
$ws.Range("E50").Value = "  -9.28%  "
$ws.Range("D51").Value = "17.69"
$ws.Range("E51").Value = "  -14.56%  "

# Row 24/25 swap: Litecoin now ranks above Polygon (source data re-sorted)
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "58.98"
$ws.Range("E24").Value = "  -15.43%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "0.439"
$ws.Range("E25").Value = "  -14.72%  "

# Row 44/45 swap: ONDO now ranks above Mantle (source data re-sorted)
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").Value = "0.944"
$ws.Range("E44").Value = "  -13.66%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.601"
$ws.Range("E45").Value = "  -13.92%  "

# Restore the default (no explicit style) look for column D so the
# persisted XML has no spurious style index, matching the original file.
$dCells.Style = "Normal"
